# Insert a new row at position 259 (this pushes the existing rows 259-316
# down to 260-317) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(259).Insert()

$ws.Cells.Item(259, 1).Value  = 5
$ws.Cells.Item(259, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(259, 3).Value  = "Maule"
$ws.Cells.Item(259, 4).Value  = 44511
$ws.Cells.Item(259, 5).Value  = 7
$ws.Cells.Item(259, 6).Value  = 100114001
$ws.Cells.Item(259, 7).Value  = "Papa"
$ws.Cells.Item(259, 8).Value  = "Rodeo"
$ws.Cells.Item(259, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(259, 10).Value = 1600
$ws.Cells.Item(259, 11).Value = 9000
$ws.Cells.Item(259, 12).Value = 9000
$ws.Cells.Item(259, 13).Value = 9000
$ws.Cells.Item(259, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(259, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(259, 16).Value = 360
$ws.Cells.Item(259, 17).Value = 25
$ws.Cells.Item(259, 18).Value = "Hortaliza"
